$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Rows 16 and 17 are brand new; give column A the same style used by the
# other index cells (bold, bordered, centered) by copying the format from
# an existing formatted cell.
$ws.Range("A15").Copy()
$ws.Range("A16:A17").PasteSpecial(-4122)
$excel.CutCopyMode = $false

# New data for rows 8 through 17 (columns A:E)
# A = index, B = name, C = from_bus, D = to_bus, E = in_service
$data = @(
    @(6,  "line7", 14, 11, $true),
    @(7,  "line8", 16, 9,  $false),
    @(8,  "extr1", 5,  12, $true),
    @(9,  "extr2", 5,  9,  $true),
    @(10, "extr3", 10, 11, $true),
    @(11, "extr4", 7,  8,  $false),
    @(12, "extr5", 9,  11, $false),
    @(13, "extr6", 7,  11, $true),
    @(14, "extr7", 5,  7,  $true),
    @(15, "extr8", 8,  5,  $false)
)

$row = 8
foreach ($rec in $data) {
    $ws.Cells.Item($row, 1).Value = $rec[0]
    $ws.Cells.Item($row, 2).Value = $rec[1]
    $ws.Cells.Item($row, 3).Value = $rec[2]
    $ws.Cells.Item($row, 4).Value = $rec[3]
    $ws.Cells.Item($row, 5).Value = $rec[4]
    $row = $row + 1
}
